$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo "Invalid usernama and password" -> "Invalid username or password"
# in both existing rows that reference it.
$ws.Range("D2").Value = "Invalid username or password"
$ws.Range("D3").Value = "Invalid username or password"

# Add a new row 4 with a fresh test case.
$ws.Range("A4").Value = "bala"
$ws.Range("B4").Value = "bala123"
$ws.Range("C4").Value = "Dutch"
$ws.Range("D4").Value = "Invalid username or password123"

# Update the active selection to D4 as in the target workbook.
$ws.Range("D4").Select()
